# hs-logger_priorities.xlsx update
# - Mark HP Mux, F250 Bridge, Nanovoltmeter and F17 Bridge as no longer
#   important (Importance 1 -> 0) now that they are complete.
# - Put transform on display / Comment points files: Difficulty 5 -> 4.
# - Add a new Todo row: "Fix datafile filenames" (Importance 1, Importance 2,
#   Difficulty 3).
# - Re-sort the whole table by Priority (column E) descending, same as the
#   worksheet's existing sort order.
# - Keep dimension / conditional-formatting ranges / selection in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing rows (by current, pre-sort position) -------------

# F250 Bridge -> importance (col B) done
$ws.Cells.Item(2, 2).Value = 0
# HP Mux -> importance (col B) done
$ws.Cells.Item(3, 2).Value = 0
# Nanovoltmeter -> importance (col B) done
$ws.Cells.Item(4, 2).Value = 0
# Put transform on display -> difficulty (col D) 5 -> 4
$ws.Cells.Item(12, 4).Value = 4
# F17 Bridge -> importance (col B) done
$ws.Cells.Item(16, 2).Value = 0
# Comment points files -> difficulty (col D) 5 -> 4
$ws.Cells.Item(17, 4).Value = 4

# --- 2. Add the new todo item as row 31 -----------------------------------

$ws.Cells.Item(31, 1).Value = "Fix datafile filenames"
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = 2
$ws.Cells.Item(31, 4).Value = 3
$ws.Cells.Item(31, 5).Formula = "=B31*C31/D31"

# --- 3. Re-sort A2:E31 by Priority (column E), descending -----------------
# (Header:=xlNo/2 -- the chosen range already excludes the header row 1,
#  so it must NOT be treated as containing a header row itself.)

$ws.Range("A2:E31").Sort($ws.Range("E2:E31"), 2, $null, $null, 1, $null, 1, 2)

# --- 4. Selection shown when the workbook was saved ------------------------

$ws.Range("C9").Select()

# --- 5. Extend the conditional-formatting ranges that covered column C/D/E
#        down to the new last row (31) ------------------------------------

$ws.Range("D2:D30").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D31"))
$ws.Range("E2:E30").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E2:E31"))
$ws.Range("C2:C30").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("C2:C31"))
